$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Update the "Approved/Rejected" column (I) for rows 3-7 to "Rejected"
# and fill in the "ReasonToReject" column (J) for rows 3-7 with "Nil".
for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = "Rejected"
    $ws.Cells.Item($r, 10).Value = "Nil"
}

# Reflect the new selection used when making the edit
$ws.Range("J3:J7").Select()
